# Collapse the paragraph describing the update command into a single
# run with new wording (no more embedded verbatim "bash -c ..." snippet).
$d = $word.ActiveDocument

$old = "When these files are placed in their respective locations, the system will run the command bash -c 'dnf update -y && dnf upgrade -y' daily."
$new = "When these files are placed in their respective locations, the system will run the update commands of various pacakge managers. Distros unsupported by a package manager will simply fail to run for that package manager.ast"

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    # Fallback: locate the paragraph by its distinctive leading text and
    # replace its whole range (covers the case where the Find above
    # doesn't match verbatim, e.g. due to smart-quote substitution).
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "When these files are placed*") {
            $start = $p.Range.Start
            $end = $p.Range.End
            $r = $d.Range($start, $end)
            $r.Text = $new
            break
        }
    }
}
